# Re-sort the worksheet tabs so that "总计" (the summary sheet) comes
# before "2021-Q1" (the quarterly detail sheet). Excel keeps each sheet's
# data/name attached as it physically moves, so this reproduces the
# commit's sheet-order swap.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$firstSheet = $wb.Worksheets.Item(1)

# Move "总计" in front of whatever sheet is currently first (the
# "2021-Q1" sheet), so the tab order becomes: 总计, 2021-Q1.
$summarySheet.Move($firstSheet)
